# Update the "想去人数" (F column) figures on the 展览 and 全部类型 sheets.
# Both sheets share the same set of rows/values that need updating.

$wb = $excel.ActiveWorkbook

$updates = @{
    2  = 8880
    3  = 8350
    5  = 168
    6  = 214
    8  = 762
    10 = 5521
    11 = 12
    17 = 170
    18 = 224
    19 = 18
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}

$wb.Save()
